$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.01

$ws.Range("H4").Value = 0

$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0.01
